$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.416.49"
$ws.Range("D3").Value = "'3.182.76"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'593.87"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'148.77"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'3.176.29"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "'6.08"
$ws.Range("E11").Value = "  +7.98%  "
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "'37.88"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "'3.721.18"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'7.37"
$ws.Range("E17").Value = "  +4.40%  "
$ws.Range("D18").Value = "'3.182.43"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "'64.198.43"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").Value = "'473.30"
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("E21").Value = "  +2.94%  "
$ws.Range("D22").Value = "'0.739"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("D23").Value = "'7.70"
$ws.Range("E23").Value = "  +3.80%  "
$ws.Range("D24").Value = "'2.45"
$ws.Range("E24").Value = "  +9.25%  "
$ws.Range("D25").Value = "'13.31"
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").Value = "'81.90"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  +8.67%  "
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "'7.26"
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  +10.58%  "
$ws.Range("D34").Value = "'28.45"
$ws.Range("E34").Value = "  +6.83%  "
$ws.Range("D35").Value = "'0.0₃0861"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'3.37"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "'467.66"
$ws.Range("E40").Value = "  +7.42%  "
$ws.Range("D41").Value = "'9.45"
$ws.Range("E41").Value = "  +8.87%  "
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("E43").Value = "  +7.64%  "
$ws.Range("D44").Value = "'0.0378"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").Value = "'2.940.59"
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("D46").Value = "'39.45"
$ws.Range("D47").Value = "'0.110"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "'132.85"
$ws.Range("E48").Value = "  +6.16%  "
$ws.Range("D50").Value = "'2.27"
$ws.Range("E50").Value = "  +6.00%  "
$ws.Range("E51").Value = "  +1.57%  "
